# Update crypto price/volume figures per the latest scrape.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.854.40'
$ws.Range("E2").Value = '  -1.54%  '
$ws.Range("D3").Value = '1.562.54'
$ws.Range("E3").Value = '  -0.06%  '
$ws.Range("E4").Value = '  -0.40%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '205.83'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.31%  '
$ws.Range("E6").Value = '  -1.54%  '
$ws.Range("E7").Value = '  -0.28%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '21.72'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -1.41%  '
$ws.Range("E9").Value = '  -0.42%  '
$ws.Range("E10").Value = '  -1.39%  '
$ws.Range("E11").Value = '  +0.26%  '
$ws.Range("D12").Value = '1.783.26'
$ws.Range("E12").Value = '  -0.51%  '
$ws.Range("D13").Value = '1.556.43'
$ws.Range("E13").Value = '  -0.47%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '3.73'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -1.87%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.513'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -0.50%  '
$ws.Range("D16").Value = '26.842.42'
$ws.Range("E16").Value = '  -1.86%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '61.14'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -3.45%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '214.06'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +1.29%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '7.35'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +1.78%  '
$ws.Range("E20").Value = '  -1.59%  '
$ws.Range("E21").Value = '  -0.04%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.12'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +0.10%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '9.18'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -2.71%  '
$ws.Range("E24").Value = '  -0.24%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '154.03'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +0.73%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '6.70'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +1.05%  '
$ws.Range("E27").Value = '  +0.20%  '
$ws.Range("E28").Value = '  -0.25%  '
$ws.Range("E29").Value = '  -1.45%  '
$ws.Range("E30").Value = '  -1.51%  '
$ws.Range("E31").Value = '  -3.40%  '
$ws.Range("E32").Value = '  -0.36%  '
$ws.Range("D33").Value = '1.400.90'
$ws.Range("E33").Value = '  +1.13%  '
$ws.Range("E34").Value = '  -0.99%  '
$ws.Range("E35").Value = '  -1.47%  '
$ws.Range("E36").Value = '  -1.22%  '
$ws.Range("E37").Value = '  -2.31%  '
$ws.Range("E38").Value = '  -0.26%  '
$ws.Range("E39").Value = '  -0.10%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.813'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -0.25%  '
$ws.Range("E41").Value = '  -0.09%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.998'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +0.43%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '5.31'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +1.12%  '
$ws.Range("E44").Value = '  -3.88%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '63.00'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -1.02%  '
$ws.Range("D47").Value = '1.696.63'
$ws.Range("E47").Value = '  -0.51%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '85.98'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +0.48%  '
$ws.Range("D49").Value = '0.0₇0987'
$ws.Range("E49").Value = '  -0.97%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0504'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +2.21%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0945'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -0.11%  '
